$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Intermediate for Mapping" ---
# Fix misspelling "Mortendad" -> "Mortandad" in column P (rows 5-8)
$wsMap = $wb.Worksheets.Item("Intermediate for Mapping")
$wsMap.Range("P5").Value = "Mortandad"
$wsMap.Range("P6").Value = "Mortandad"
$wsMap.Range("P7").Value = "Mortandad"
$wsMap.Range("P8").Value = "Mortandad"

# --- Sheet 2: "Intermediate Exhibit" ---
$wsEx = $wb.Worksheets.Item("Intermediate Exhibit")

# Fix misspelling in the "Mortendad Canyon" section header
$wsEx.Range("A7").Value = "Mortandad Canyon"

# Widen the Longitude column slightly
$wsEx.Columns.Item(3).ColumnWidth = 13

# Split the combined "Los Alamos and Pajarito Canyons" section into two:
# row 12 keeps the Los Alamos wells and becomes "Los Alamos Canyon";
# a brand-new "Pajarito Canyon" section header is inserted as row 19,
# pushing the former row 19 (well 03-B-10) down to row 20.
$wsEx.Range("A12").Value = "Los Alamos Canyon"

# Duplicate the current row 19 (03-B-10 data) down into row 20 first so we
# don't lose it once row 19 is turned into the new section header.
$wsEx.Range("A19:H19").Copy()
$wsEx.Range("A20:H20").PasteSpecial(-4104)
$wsEx.Range("A20:H20").PasteSpecial(-4122)

# Merge row 19 across A:H *before* pasting the section-header formatting so
# the engine reuses the existing "section header" style (like rows 3/7/12)
# instead of synthesizing a brand-new split-border style.
$wsEx.Range("A19:H19").Merge()
$wsEx.Range("A12:H12").Copy()
$wsEx.Range("A19:H19").PasteSpecial(-4122)
$wsEx.Range("A19").Value = "Pajarito Canyon"

$excel.CutCopyMode = 0
